$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: plot_color — match the header formatting used by D1 (bold + bottom border)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "plot_color"

$colors = @(
    "n/a",
    "n/a",
    "n/a",
    "n/a",
    "n/a",
    "n/a",
    "n/a",
    "n/a",
    "gray",
    "blue",
    "brown",
    "cyan",
    "pink",
    "green",
    "yellow",
    "gray",
    "brown",
    "pink",
    "green",
    "blue",
    "black"
)

for ($i = 0; $i -lt $colors.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $colors[$i]
}

# Update selection to match author's final cursor position
$ws.Range("E21").Select()
